# Apply cell value updates per the commit diff for cryptos.xlsx (sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.846.65"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.629.35"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0881"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "1.860.25"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "1.630.52"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "27.863.00"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.411.36"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").Value = "  +2.49%  "
$ws.Range("E36").Value = "  -3.05%  "
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0170"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.553"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("E41").Value = "  -2.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("E46").Value = "  -3.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0502"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("E51").Value = "  -0.27%  "
